$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Business Cards")

# Data for rows 2-7 (columns A-H)
$data = @(
    @(1, "Nehaal Fakih", "SPRINGBOARD", "", "+91-9004488330", "nehaal@91springboard.com", "91springboard.com", "Plot No. D-5 Road No. 20, Marol MIDC, Andheri East"),
    @(2, "Tel", "SRIMATHA", "", "(08814)224530, 224796 (0), 94401 80153", "Ivrcocoanutmerchants@gmail.com", "M.Narayana", ""),
    @(3, "CA", "Near Hotel", "", "91 9730704929, 91 9421833600", "capushkarsolanki@gmail.com", "B.Com", ""),
    @(4, "Bhavesh", "Shreenath", "Sales Corporation", "93253 13803, 99675 53803", "", "", "Off Link Road, Mahavir Nagar;"),
    @(5, "Plywood", "Laminates", "", "9370762286, 9373962286, 0253-2462286", "", "", "Nashik Road"),
    @(6, "MANAVTA", "Gram", "", "(0) 6501982, 4618756, 4618738", "", "", "")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowData = $data[$i]

    if ($row -ne 2) {
        # Copy style from the existing data row 2 template (col A style s=2, cols B-H style s=3)
        $ws.Cells.Item(2, 1).Copy() | Out-Null
        $ws.Cells.Item($row, 1).PasteSpecial(-4122) | Out-Null
        for ($c = 2; $c -le 8; $c++) {
            $ws.Cells.Item(2, 2).Copy() | Out-Null
            $ws.Cells.Item($row, $c).PasteSpecial(-4122) | Out-Null
        }
    }

    for ($c = 1; $c -le 8; $c++) {
        $ws.Cells.Item($row, $c).Value = $rowData[$c - 1]
    }
}

$excel.CutCopyMode = $false

# Update autofilter range to cover the new extent
$ws.AutoFilterMode = $false
$ws.Range("A1:H7").AutoFilter() | Out-Null

# Update the hidden _FilterDatabase defined name to match the new range
foreach ($dn in $wb.Names) {
    if ($dn.Name -like "*_FilterDatabase*") {
        $dn.RefersTo = "='Business Cards'!`$A`$1:`$H`$7"
    }
}

